$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price + 1h volume change).
# D-column "prices" are stored as literal text (they use "." as a thousands
# separator in several rows, e.g. "25.786.75"), so force text formatting
# before writing to stop Excel from auto-converting them to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.786.75'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.745.43'
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.40'

$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5084'
$ws.Range('E7').Value = '  +3.19%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.36'
$ws.Range('E8').Value = '  -2.79%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2664'
$ws.Range('E9').Value = '  +4.53%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06160'
$ws.Range('E10').Value = '  +2.47%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.753.75'
$ws.Range('E11').Value = '  +0.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06933'
$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.31'
$ws.Range('E13').Value = '  +3.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6224'
$ws.Range('E14').Value = '  +9.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.465'
$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.56'
$ws.Range('E16').Value = '  +1.52%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.817.98'
$ws.Range('E19').Value = '  +0.38%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.55'
$ws.Range('E20').Value = '  +2.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006615'
$ws.Range('E21').Value = '  +1.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.971.24'
$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.049'
$ws.Range('E23').Value = '  +1.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.235'
$ws.Range('E24').Value = '  +4.02%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.132'
$ws.Range('E25').Value = '  +1.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.16'
$ws.Range('E26').Value = '  -0.61%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.467'
$ws.Range('E27').Value = '  -0.52%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.02'
$ws.Range('E28').Value = '  +2.54%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.768'
$ws.Range('E29').Value = '  -2.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.35'
$ws.Range('E30').Value = '  +0.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08182'
$ws.Range('E31').Value = '  +2.77%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.669'
$ws.Range('E32').Value = '  -1.98%  '

$ws.Range('E33').Value = '  -0.58%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04371'
$ws.Range('E34').Value = '  -0.27%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.647'
$ws.Range('E35').Value = '  +1.51%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9928'
$ws.Range('E36').Value = '  +1.41%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5965'
$ws.Range('E37').Value = '  +0.10%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.612'
$ws.Range('E38').Value = '  -2.05%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01553'
$ws.Range('E39').Value = '  +2.85%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.911'
$ws.Range('E40').Value = '  -0.47%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.47'
$ws.Range('E42').Value = '  -0.26%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3812'
$ws.Range('E43').Value = '  +1.61%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7438'
$ws.Range('E44').Value = '  -1.15%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.870'
$ws.Range('E45').Value = '  -5.32%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05487'
$ws.Range('E46').Value = '  +4.93%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1091'
$ws.Range('E47').Value = '  +2.31%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.924'
$ws.Range('E48').Value = '  +2.33%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.97'
$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.37'
$ws.Range('E50').Value = '  +0.57%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.30%  '
